# "changed pictures and load strategy"
# The workbook's question_img_name column (C) used to hold long,
# descriptive GIF file names (e.g. "GIF3_Q2_MysteryDrink.gif"). The new
# loading strategy simply references the GIFs by their numeric id
# ("3.gif", "4.gif", ...), so update each of those five cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "3.gif"
$ws.Range("C3").Value = "4.gif"
$ws.Range("C4").Value = "5.gif"
$ws.Range("C5").Value = "6.gif"
$ws.Range("C6").Value = "7.gif"

# Leave the cursor where the author ended up after making the edits.
$ws.Range("C7").Select()
